$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.169.06'
$ws.Range("E2").Value = '  -1.04%  '
$ws.Range("D3").Value = '1.572.73'
$ws.Range("E3").Value = '  -0.38%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("E6").Value = '  -1.37%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.26'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.40%  '
$ws.Range("E9").Value = '  -0.64%  '
$ws.Range("E10").Value = '  +0.17%  '
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("D12").Value = '1.796.36'
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("D13").Value = '1.571.72'
$ws.Range("E13").Value = '  -0.38%  '
$ws.Range("E14").Value = '  -1.20%  '
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").Value = '27.186.87'
$ws.Range("E16").Value = '  -1.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.23'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.22'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.14%  '
$ws.Range("E19").Value = '  +1.11%  '
$ws.Range("E20").Value = '  -0.59%  '
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("E23").Value = '  -3.63%  '
$ws.Range("E24").Value = '  +0.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.67'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.17%  '
$ws.Range("E26").Value = '  -3.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.93'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("E29").Value = '  -1.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.55%  '
$ws.Range("E31").Value = '  -1.62%  '
$ws.Range("E32").Value = '  -1.18%  '
$ws.Range("D33").Value = '1.409.97'
$ws.Range("E33").Value = '  +3.56%  '
$ws.Range("E34").Value = '  -1.06%  '
$ws.Range("E35").Value = '  +2.22%  '
$ws.Range("E36").Value = '  -0.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.944'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.31%  '
$ws.Range("E38").Value = '  -1.73%  '
$ws.Range("E39").Value = '  -0.49%  '
$ws.Range("E40").Value = '  -2.66%  '
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.997'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.79%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.83'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.39'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.20'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.69'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.55%  '
$ws.Range("D47").Value = '1.708.70'
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.73'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.78%  '
$ws.Range("D49").Value = '0.0₇0983'
$ws.Range("E49").Value = '  -1.45%  '
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("E51").Value = '  -0.04%  '
